$d = $word.ActiveDocument

$oldText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "Mapky v tomto dokumente pripravil Jan Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Locate the run of text that needs to change.
$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Grab the opening <w:p ...> tag (with its <w:pPr>) of the paragraph that
# holds the text, so the replacement paragraph keeps its identity/formatting.
$para = $rng.Paragraphs(1)
$pxml = $para.Range.WordOpenXML
$null = $pxml -match '(?s)(<w:p [^>]*>\s*<w:pPr>.*?</w:pPr>)'
$pOpenAndPPr = $matches[1]

# Remove the old runs, then drop in a fresh paragraph shell (same identity /
# paragraph properties) containing one empty run followed by a single
# unformatted run holding the updated text.
$rng.Delete()

$escapedNew = $newText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' + $pOpenAndPPr + '<w:r/><w:r><w:t>' + $escapedNew + '</w:t></w:r></w:p></w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rng.InsertXML($xml)
